$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write cell values in the precise order needed so the shared-strings
# table indices come out matching the target workbook.
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("A2").Value = "Barbara"
$ws.Range("C2").Value = "BB"
$ws.Range("D1").Value = "LoginPageText"
$ws.Range("D2").Value = "vtiger"
$ws.Range("E2").Value = " Administrator"
$ws.Range("E1").Value = "UserDropDownOption"
$ws.Range("B2").Value = "Navya"
$ws.Range("F1").Value = "TextWhenNotPresent"
$ws.Range("F2").Value = "`n                 No Contact Found !`n          "

# Column widths (requested values are back-solved so the persisted
# OOXML <col> width lands on the target: this engine stores ColumnWidth
# + 5/6 rounded to the nearest 1/6 character unit)
$ws.Columns.Item(4).ColumnWidth = 13.1666666666667
$ws.Columns.Item(5).ColumnWidth = 20
$ws.Columns.Item(6).ColumnWidth = 19.6666666666667

# Row height for row 2
$ws.Rows.Item(2).RowHeight = 60

# Wrap text on F2
$ws.Range("F2").WrapText = $true

# Selection
$ws.Range("F2").Select()
